$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 31250318
$ws.Range("I6").Value = 38461812
$ws.Range("K6").Value = 115385436
$ws.Range("M6").Value = -115385324
$ws.Range("H8").Value = 1319
$ws.Range("I8").Value = 1403.75
$ws.Range("K8").Value = 4211.25
$ws.Range("M8").Value = -4072.25
$ws.Range("H40").Value = 3196.125
$ws.Range("J40").Value = 3175.4167
$ws.Range("L40").Value = 3175.4167
$ws.Range("N40").Value = -3525.4167
$ws.Range("H51").Value = 3399.8
$ws.Range("J51").Value = 3666.3333
$ws.Range("L51").Value = 3666.3333
$ws.Range("N51").Value = -4634.3333
$ws.Range("H106").Value = 8328.956
$ws.Range("J106").Value = 9620.143
$ws.Range("L106").Value = 9620.143
$ws.Range("N106").Value = -10882.143
$ws.Range("H129").Value = 2454.2
$ws.Range("J129").Value = 2981.8333
$ws.Range("L129").Value = 8945.499899999999
$ws.Range("N129").Value = -18945.4999
$ws.Range("H132").Value = 2149.182
$ws.Range("I132").Value = 1485.4667
$ws.Range("K132").Value = 4456.4001
$ws.Range("M132").Value = -1926.4001
$ws.Range("H137").Value = 2529.4814
$ws.Range("I137").Value = 2253.3333
$ws.Range("K137").Value = 6759.999899999999
$ws.Range("M137").Value = -4209.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 3742
$ws.Range("I36").Value = 3742
$ws.Range("K36").Value = 3742
$ws.Range("M36").Value = -3396
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null
$ws.Range("H124").Value = 52806.668
$ws.Range("J124").Value = 52806.668
$ws.Range("L124").Value = 52806.668
$ws.Range("N124").Value = -62626.668
$ws.Range("H125").Value = 58718.332
$ws.Range("J125").Value = 58718.332
$ws.Range("L125").Value = 58718.332
$ws.Range("N125").Value = -68558.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1838.9149
$ws.Range("I134").Value = 1706.975
$ws.Range("K134").Value = 5120.924999999999
$ws.Range("M134").Value = -2585.924999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4331.709
$ws.Range("I31").Value = 2401.6667
$ws.Range("J31").Value = 5055.475
$ws.Range("K31").Value = 2401.6667
$ws.Range("L31").Value = 5055.475
$ws.Range("M31").Value = -2106.6667
$ws.Range("N31").Value = -5645.475
$ws.Range("H34").Value = 4331.709
$ws.Range("I34").Value = 2401.6667
$ws.Range("J34").Value = 5055.475
$ws.Range("K34").Value = 2401.6667
$ws.Range("L34").Value = 5055.475
$ws.Range("M34").Value = -2199.6667
$ws.Range("N34").Value = -5459.475
$ws.Range("H58").Value = 2582.375
$ws.Range("I58").Value = 2983.0715
$ws.Range("J58").Value = 2021.4
$ws.Range("K58").Value = 2983.0715
$ws.Range("L58").Value = 2021.4
$ws.Range("M58").Value = -2780.0715
$ws.Range("N58").Value = -2427.4
$ws.Range("H74").Value = 45828.5
$ws.Range("J74").Value = 46104.668
$ws.Range("L74").Value = 46104.668
$ws.Range("N74").Value = -47852.668
$ws.Range("H77").Value = 45828.5
$ws.Range("J77").Value = 46104.668
$ws.Range("L77").Value = 138314.004
$ws.Range("N77").Value = -147050.004
$ws.Range("H133").Value = 60326
$ws.Range("J133").Value = 60326
$ws.Range("L133").Value = 60326
$ws.Range("N133").Value = -65386
$ws.Range("H136").Value = 2582.375
$ws.Range("I136").Value = 2983.0715
$ws.Range("J136").Value = 2021.4
$ws.Range("K136").Value = 8949.2145
$ws.Range("L136").Value = 6064.200000000001
$ws.Range("M136").Value = -6399.2145
$ws.Range("N136").Value = -11164.2
$ws.Range("H141").Value = 198333
$ws.Range("J141").Value = 198333
$ws.Range("L141").Value = 198333
$ws.Range("N141").Value = -208693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 95.1579
$ws.Range("J12").Value = 81.09090999999999
$ws.Range("L12").Value = 243.27273
$ws.Range("N12").Value = -589.2727299999999
$ws.Range("H109").Value = 7634
$ws.Range("I109").Value = 662.3333
$ws.Range("K109").Value = 1986.9999
$ws.Range("M109").Value = -946.9999
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = $null
$ws.Range("N116").Value = $null
$ws.Range("H131").Value = 2072.8696
$ws.Range("I131").Value = 1157
$ws.Range("J131").Value = 2265.6843
$ws.Range("K131").Value = 3471
$ws.Range("L131").Value = 6797.0529
$ws.Range("M131").Value = 1569
$ws.Range("N131").Value = -16877.0529

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 22997
$ws.Range("J15").Value = 22997
$ws.Range("L15").Value = 22997
$ws.Range("N15").Value = -23573
$ws.Range("H81").Value = 22997
$ws.Range("J81").Value = 22997
$ws.Range("L81").Value = 22997
$ws.Range("N81").Value = -24993
$ws.Range("H84").Value = 22997
$ws.Range("J84").Value = 22997
$ws.Range("L84").Value = 68991
$ws.Range("N84").Value = -78975
$ws.Range("H132").Value = 2664.2307
$ws.Range("I132").Value = 1830.6316
$ws.Range("J132").Value = 4926.857
$ws.Range("K132").Value = 5491.8948
$ws.Range("L132").Value = 14780.571
$ws.Range("M132").Value = -2961.8948
$ws.Range("N132").Value = -19840.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3085.4707
$ws.Range("J46").Value = 3316.8667
$ws.Range("L46").Value = 3316.8667
$ws.Range("N46").Value = -3692.8667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -1336
$ws.Range("H81").Value = 1448
$ws.Range("I81").Value = 1296.5
$ws.Range("K81").Value = 2593
$ws.Range("M81").Value = -1532
$ws.Range("H84").Value = 1448
$ws.Range("I84").Value = 1296.5
$ws.Range("K84").Value = 12965
$ws.Range("M84").Value = -7661
$ws.Range("H122").Value = 7557.222
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 8189.375
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 24568.125
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -29468.125
$ws.Range("H126").Value = 3008.1614
$ws.Range("I126").Value = 2783.5417
$ws.Range("K126").Value = 8350.625100000001
$ws.Range("M126").Value = -5880.625100000001
$ws.Range("H136").Value = 1922.9131
$ws.Range("I136").Value = 1514
$ws.Range("K136").Value = 4542
$ws.Range("M136").Value = -1992
